$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 1017
$ws.Range("I21").Value = 1017
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 1017
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -549
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 1017
$ws.Range("I23").Value = 1017
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 1017
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -783
$ws.Range("N23").ClearContents()
$ws.Range("H62").Value = 37042252
$ws.Range("I62").Value = 55560580
$ws.Range("J62").Value = 5588.8887
$ws.Range("K62").Value = 55560580
$ws.Range("L62").Value = 5588.8887
$ws.Range("M62").Value = -55559956
$ws.Range("N62").Value = -6836.8887
$ws.Range("H65").Value = 37042252
$ws.Range("I65").Value = 55560580
$ws.Range("J65").Value = 5588.8887
$ws.Range("K65").Value = 277802900
$ws.Range("L65").Value = 27944.4435
$ws.Range("M65").Value = -277799780
$ws.Range("N65").Value = -34184.4435
$ws.Range("H82").Value = 197
$ws.Range("I82").Value = 197
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 591
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -185
$ws.Range("H85").Value = 197
$ws.Range("I85").Value = 197
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 591
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = 813
$ws.Range("H98").Value = 32931.1
$ws.Range("I98").Value = 33154.53
$ws.Range("J98").Value = 31665
$ws.Range("K98").Value = 33154.53
$ws.Range("L98").Value = 31665
$ws.Range("M98").Value = -31656.53
$ws.Range("N98").Value = -34661
$ws.Range("H122").Value = 32931.1
$ws.Range("I122").Value = 33154.53
$ws.Range("J122").Value = 31665
$ws.Range("K122").Value = 99463.59
$ws.Range("L122").Value = 94995
$ws.Range("M122").Value = -97013.59
$ws.Range("N122").Value = -99895
$ws.Range("H125").Value = 2600
$ws.Range("I125").Value = 2000
$ws.Range("J125").Value = 3500
$ws.Range("K125").Value = 18000
$ws.Range("L125").Value = 31500
$ws.Range("M125").Value = -15540
$ws.Range("N125").Value = -36420

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 93666
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 93666
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 93666
$ws.Range("N44").Value = -94642
$ws.Range("H88").Value = 142860510
$ws.Range("I88").Value = 2865.3333
$ws.Range("J88").Value = 250003740
$ws.Range("K88").Value = 2865.3333
$ws.Range("L88").Value = 250003740
$ws.Range("M88").Value = -2459.3333
$ws.Range("N88").Value = -250004552
$ws.Range("H91").Value = 142860510
$ws.Range("I91").Value = 2865.3333
$ws.Range("J91").Value = 250003740
$ws.Range("K91").Value = 2865.3333
$ws.Range("L91").Value = 250003740
$ws.Range("M91").Value = -1461.3333
$ws.Range("N91").Value = -250006548

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2113.0417
$ws.Range("I20").Value = 1141.9333
$ws.Range("J20").Value = 3731.5557
$ws.Range("K20").Value = 1141.9333
$ws.Range("L20").Value = 3731.5557
$ws.Range("M20").Value = -894.9332999999999
$ws.Range("N20").Value = -4225.5557
$ws.Range("H56").Value = 19750
$ws.Range("I56").Value = 25000
$ws.Range("J56").Value = 14500
$ws.Range("K56").Value = 25000
$ws.Range("L56").Value = 14500
$ws.Range("M56").Value = -24261
$ws.Range("N56").Value = -15978
$ws.Range("H86").Value = 6160.6
$ws.Range("I86").Value = 4983.773
$ws.Range("J86").Value = 9396.875
$ws.Range("K86").Value = 4983.773
$ws.Range("L86").Value = 9396.875
$ws.Range("M86").Value = -3860.773
$ws.Range("N86").Value = -11642.875
$ws.Range("H89").Value = 6160.6
$ws.Range("I89").Value = 4983.773
$ws.Range("J89").Value = 9396.875
$ws.Range("K89").Value = 24918.865
$ws.Range("L89").Value = 46984.375
$ws.Range("M89").Value = -19302.865
$ws.Range("N89").Value = -58216.375
$ws.Range("H134").Value = 3377.9722
$ws.Range("I134").Value = 2949.1924
$ws.Range("J134").Value = 4492.8
$ws.Range("K134").Value = 8847.5772
$ws.Range("L134").Value = 13478.4
$ws.Range("M134").Value = -6312.5772
$ws.Range("N134").Value = -18548.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2005.3793
$ws.Range("I58").Value = 1100.3478
$ws.Range("J58").Value = 5474.6665
$ws.Range("K58").Value = 1100.3478
$ws.Range("L58").Value = 5474.6665
$ws.Range("M58").Value = -897.3478
$ws.Range("N58").Value = -5880.6665
$ws.Range("H70").Value = 59998.5
$ws.Range("I70").Value = 79998
$ws.Range("J70").Value = 39999
$ws.Range("K70").Value = 79998
$ws.Range("L70").Value = 39999
$ws.Range("M70").Value = -79683
$ws.Range("N70").Value = -40629
$ws.Range("H73").Value = 59998.5
$ws.Range("I73").Value = 79998
$ws.Range("J73").Value = 39999
$ws.Range("K73").Value = 79998
$ws.Range("L73").Value = 39999
$ws.Range("M73").Value = -78906
$ws.Range("N73").Value = -42183
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H134").Value = 4610.926
$ws.Range("I134").Value = 2212.5557
$ws.Range("J134").Value = 9407.666999999999
$ws.Range("K134").Value = 6637.6671
$ws.Range("L134").Value = 28223.001
$ws.Range("M134").Value = -4102.6671
$ws.Range("N134").Value = -33293.001
$ws.Range("H136").Value = 2005.3793
$ws.Range("I136").Value = 1100.3478
$ws.Range("J136").Value = 5474.6665
$ws.Range("K136").Value = 3301.0434
$ws.Range("L136").Value = 16423.9995
$ws.Range("M136").Value = -751.0434
$ws.Range("N136").Value = -21523.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 386421.06
$ws.Range("I5").Value = 1396
$ws.Range("J5").Value = 590257.9
$ws.Range("K5").Value = 4188
$ws.Range("L5").Value = 1770773.7
$ws.Range("M5").Value = -4076
$ws.Range("N5").Value = -1770997.7
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H52").Value = 25565.5
$ws.Range("I52").Value = 499.5
$ws.Range("J52").Value = 50631.5
$ws.Range("K52").Value = 1498.5
$ws.Range("L52").Value = 151894.5
$ws.Range("M52").Value = -1232.5
$ws.Range("N52").Value = -152426.5
$ws.Range("H135").Value = 386421.06
$ws.Range("I135").Value = 1396
$ws.Range("J135").Value = 590257.9
$ws.Range("K135").Value = 12564
$ws.Range("L135").Value = 5312321.100000001
$ws.Range("M135").Value = -10029
$ws.Range("N135").Value = -5317391.100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2500
$ws.Range("I80").Value = 2000
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 2000
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -1002
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 2500
$ws.Range("I83").Value = 2000
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 10000
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -5008
$ws.Range("N83").Value = -24984
$ws.Range("H123").Value = 31583.334
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 31583.334
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 31583.334
$ws.Range("N123").Value = -36483.334
$ws.Range("H132").Value = 3439.0557
$ws.Range("I132").Value = 3768.3845
$ws.Range("J132").Value = 2582.8
$ws.Range("K132").Value = 11305.1535
$ws.Range("L132").Value = 7748.400000000001
$ws.Range("M132").Value = -8775.1535
$ws.Range("N132").Value = -12808.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 23239.13
$ws.Range("I7").Value = 31266.6
$ws.Range("J7").Value = 8187.625
$ws.Range("K7").Value = 31266.6
$ws.Range("L7").Value = 8187.625
$ws.Range("M7").Value = -31154.6
$ws.Range("N7").Value = -8411.625
$ws.Range("H16").Value = 3478.4
$ws.Range("I16").Value = 3135.0527
$ws.Range("J16").Value = 10002
$ws.Range("K16").Value = 3135.0527
$ws.Range("L16").Value = 10002
$ws.Range("M16").Value = -2965.0527
$ws.Range("N16").Value = -10342
$ws.Range("H38").Value = 41660.5
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 41660.5
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 41660.5
$ws.Range("N38").Value = -42480.5
$ws.Range("M38").ClearContents()
$ws.Range("H61").Value = 7675.24
$ws.Range("I61").Value = 7111.4375
$ws.Range("J61").Value = 8677.556
$ws.Range("K61").Value = 7111.4375
$ws.Range("L61").Value = 8677.556
$ws.Range("M61").Value = -6909.4375
$ws.Range("N61").Value = -9081.556
$ws.Range("H68").Value = 5939.25
$ws.Range("I68").Value = 2310.5
$ws.Range("J68").Value = 7148.8335
$ws.Range("K68").Value = 2310.5
$ws.Range("L68").Value = 7148.8335
$ws.Range("M68").Value = -1561.5
$ws.Range("N68").Value = -8646.833500000001
$ws.Range("H71").Value = 5939.25
$ws.Range("I71").Value = 2310.5
$ws.Range("J71").Value = 7148.8335
$ws.Range("K71").Value = 11552.5
$ws.Range("L71").Value = 35744.1675
$ws.Range("M71").Value = -7808.5
$ws.Range("N71").Value = -43232.1675
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("M90").ClearContents()
$ws.Range("H113").Value = 7675.24
$ws.Range("I113").Value = 7111.4375
$ws.Range("J113").Value = 8677.556
$ws.Range("K113").Value = 7111.4375
$ws.Range("L113").Value = 8677.556
$ws.Range("M113").Value = -4941.4375
$ws.Range("N113").Value = -13017.556
$ws.Range("H122").Value = 5419.3687
$ws.Range("I122").Value = 4077.2727
$ws.Range("J122").Value = 7264.75
$ws.Range("K122").Value = 12231.8181
$ws.Range("L122").Value = 21794.25
$ws.Range("M122").Value = -9781.8181
$ws.Range("N122").Value = -26694.25
$ws.Range("H126").Value = 23239.13
$ws.Range("I126").Value = 31266.6
$ws.Range("J126").Value = 8187.625
$ws.Range("K126").Value = 93799.79999999999
$ws.Range("L126").Value = 24562.875
$ws.Range("M126").Value = -91329.79999999999
$ws.Range("N126").Value = -29502.875
$ws.Range("H136").Value = 3744.0425
$ws.Range("I136").Value = 1323.6129
$ws.Range("J136").Value = 8433.625
$ws.Range("K136").Value = 3970.8387
$ws.Range("L136").Value = 25300.875
$ws.Range("M136").Value = -1420.8387
$ws.Range("N136").Value = -30400.875

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 7755.9575
$ws.Range("I132").Value = 8105.4
$ws.Range("J132").Value = 6554.75
$ws.Range("K132").Value = 24316.2
$ws.Range("L132").Value = 19664.25
$ws.Range("M132").Value = -21786.2
$ws.Range("N132").Value = -24724.25
